$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "to" (wallet address) column for both rows with new values,
# replacing the old placeholder addresses.
$ws.Range("G2").Value = "0x7CE2cfC2b3c838150d110cF78d156d96674afe54"
$ws.Range("G3").Value = "0x04628CC763C09e41aDC9C6b6F28ED7d6c35d7d42"

# Update the active selection to E11 (was E9)
$ws.Range("E11").Select()
